# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the "_old" / "_new" column-header suffixes to the respective
# format-version suffixes ("_FV2410" / "_FV2504"), wraps the sheet's data
# range in a native Excel Table (so the header row is discoverable /
# filterable), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# 1) Rename the header cells (row 1) — "_old" -> "_FV2410", "_new" -> "_FV2504".
#    "diff" (column K / index 11) is left untouched.
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the used range into a proper Excel Table covering the header row
#    and all data rows (A1:U62), named "Table1".
$dataRange = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (split under row 1, top-left cell of the
#    scrollable area is A2).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
